# Commit: "Added other modules in Ruby to selendroid"
#
# The VT200-0578 test case (row 5 of the TestCases sheet) had its
# "Steps" (G5) and "Expected Behaviour" (H5) columns updated:
#   - G5: the two TakeScreenshot(...) calls flanking validate4/validate6
#         were removed (screenshots are no longer taken at those points).
#   - H5: validate4 / validate6 switched from comparing a screenshot to
#         checking whether the tabbar icon is displayed / not displayed.
# The active selection on the sheet also moved from G2 to E1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

$ws.Range("G5").Value = "wait(5);
validate1;
link_Click(tabbar_test_link);
validate2;
SelectTestToRun(VT200_0576_string);
ClickRunTest(runtest_top_xpath);
validate3;
ClickRunTest(runtest_bottom_xpath);
wait(5);
validate4;
SwitchApp(NATIVE_APP);
wait(2);
ClickNativeIcon(VT200_0576_mainpage_xpath);
wait(2);
SwitchApp(WEBVIEW);
SelectTestToRun(VT200_0578_string);
ClickRunTest(runtest_top_xpath);
validate5;
ClickRunTest(runtest_bottom_xpath);
wait(5);
validate6;"

$ws.Range("H5").Value = "validate1
{
validate_PageTitle=Compliance JS specs
};
validate2
{
validate_PageTitle=Native Tabbar JS Test
};
validate3
{
validate_Text_Exists=VT200-0576
};
validate4
{
validate_isIconDisplayed=tabbar_xpath,true
};
validate5
{
validate_Text_Exists=VT200-0578
};
validate6
{
validate_isIconDisplayed=tabbar_xpath,false
};"

# Move the active selection from G2 to E1, matching the saved sheet view.
[void]$ws.Range("E1").Select()
